# Scheduled market-data refresh: update cached Leve profit calculations
# across sheets (mirrors external data-source values).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 195.125
$ws.Range("I4").Value = 195.125
$ws.Range("K4").Value = 195.125
$ws.Range("M4").Value = -81.125
$ws.Range("H98").Value = 1139.0769
$ws.Range("I98").Value = 936.7778
$ws.Range("K98").Value = 936.7778
$ws.Range("M98").Value = 561.2222
$ws.Range("H99").Value = 628.9
$ws.Range("I99").Value = 365.55554
$ws.Range("J99").Value = 2999
$ws.Range("K99").Value = 1096.66662
$ws.Range("L99").Value = 8997
$ws.Range("M99").Value = 401.33338
$ws.Range("N99").Value = -11993
$ws.Range("H111").Value = 2564.8696
$ws.Range("I111").Value = 2934.8
$ws.Range("K111").Value = 8804.400000000001
$ws.Range("M111").Value = -5737.400000000001
$ws.Range("H122").Value = 1139.0769
$ws.Range("I122").Value = 936.7778
$ws.Range("K122").Value = 2810.3334
$ws.Range("M122").Value = -360.3334
$ws.Range("H135").Value = 2666.423
$ws.Range("J135").Value = 2926.5
$ws.Range("L135").Value = 26338.5
$ws.Range("N135").Value = -31408.5
$ws.Range("H136").Value = 186857
$ws.Range("J136").Value = 186857
$ws.Range("L136").Value = 186857
$ws.Range("N136").Value = -197057
$ws.Range("H141").Value = 2600
$ws.Range("I141").Value = 2600
$ws.Range("K141").Value = 7800
$ws.Range("M141").Value = -2620

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 36454.86
$ws.Range("I32").Value = 39722.16
$ws.Range("K32").Value = 39722.16
$ws.Range("M32").Value = -39435.16
$ws.Range("H133").Value = 119999
$ws.Range("J133").Value = 119999
$ws.Range("L133").Value = 119999
$ws.Range("N133").Value = -125059

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1410.2941
$ws.Range("I107").Value = 1297.3846
$ws.Range("K107").Value = 1297.3846
$ws.Range("M107").Value = 622.6153999999999
$ws.Range("H134").Value = 3767.78
$ws.Range("I134").Value = 2645.9
$ws.Range("J134").Value = 8255.299999999999
$ws.Range("K134").Value = 7937.700000000001
$ws.Range("L134").Value = 24765.9
$ws.Range("M134").Value = -5402.700000000001
$ws.Range("N134").Value = -29835.9

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3638.125
$ws.Range("I134").Value = 2648.4092
$ws.Range("K134").Value = 7945.2276
$ws.Range("M134").Value = -5410.2276

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2689075.5
$ws.Range("I4").Value = 3953113
$ws.Range("K4").Value = 11859339
$ws.Range("M4").Value = -11859227
$ws.Range("H7").Value = 206.4
$ws.Range("I7").Value = 181.83333
$ws.Range("J7").Value = 243.25
$ws.Range("K7").Value = 545.49999
$ws.Range("L7").Value = 729.75
$ws.Range("M7").Value = -433.49999
$ws.Range("N7").Value = -953.75
$ws.Range("H12").Value = 359.5
$ws.Range("J12").Value = 359.5
$ws.Range("L12").Value = 1078.5
$ws.Range("N12").Value = -1424.5
$ws.Range("H34").Value = 4695.4165
$ws.Range("I34").Value = 57.166668
$ws.Range("K34").Value = 171.500004
$ws.Range("M34").Value = -87.50000399999999
$ws.Range("H39").Value = 4707.067
$ws.Range("J39").Value = 5388.923
$ws.Range("L39").Value = 16166.769
$ws.Range("N39").Value = -16754.769
$ws.Range("H46").Value = 278.55554
$ws.Range("I46").Value = 29.571428
$ws.Range("K46").Value = 88.71428400000001
$ws.Range("M46").Value = 2.285715999999994
$ws.Range("H55").Value = 2243.077
$ws.Range("I55").Value = 468.75
$ws.Range("J55").Value = 3031.6667
$ws.Range("K55").Value = 1406.25
$ws.Range("L55").Value = 9095.000100000001
$ws.Range("M55").Value = -1229.25
$ws.Range("N55").Value = -9449.000100000001
$ws.Range("H80").Value = 2979.6
$ws.Range("I80").Value = 2500
$ws.Range("J80").Value = 3099.5
$ws.Range("K80").Value = 7500
$ws.Range("L80").Value = 9298.5
$ws.Range("M80").Value = -6564
$ws.Range("N80").Value = -11170.5
$ws.Range("H83").Value = 2979.6
$ws.Range("I83").Value = 2500
$ws.Range("J83").Value = 3099.5
$ws.Range("K83").Value = 22500
$ws.Range("L83").Value = 27895.5
$ws.Range("M83").Value = -17820
$ws.Range("N83").Value = -37255.5
$ws.Range("H104").Value = 2500
$ws.Range("I104").Value = 2000
$ws.Range("J104").Value = 3000
$ws.Range("K104").Value = 6000
$ws.Range("L104").Value = 9000
$ws.Range("M104").Value = -3379
$ws.Range("N104").Value = -14242
$ws.Range("H105").Value = 7000
$ws.Range("I105").Value = 7000
$ws.Range("K105").Value = 21000
$ws.Range("M105").Value = -18379
$ws.Range("H107").Value = 1701
$ws.Range("I107").Value = 210
$ws.Range("K107").Value = 630
$ws.Range("M107").Value = 1290
$ws.Range("H108").Value = 750
$ws.Range("I108").Value = 750
$ws.Range("K108").Value = 2250
$ws.Range("M108").Value = 630
$ws.Range("H122").Value = 1177.5
$ws.Range("J122").Value = 1173
$ws.Range("L122").Value = 10557
$ws.Range("N122").Value = -15457
$ws.Range("H131").Value = 14496928
$ws.Range("J131").Value = 5418.4707
$ws.Range("L131").Value = 16255.4121
$ws.Range("N131").Value = -26335.4121

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5901.6665
$ws.Range("I70").Value = 5506.2
$ws.Range("J70").Value = 6396
$ws.Range("K70").Value = 5506.2
$ws.Range("L70").Value = 6396
$ws.Range("M70").Value = -5236.2
$ws.Range("N70").Value = -6936
$ws.Range("H73").Value = 5901.6665
$ws.Range("I73").Value = 5506.2
$ws.Range("J73").Value = 6396
$ws.Range("K73").Value = 5506.2
$ws.Range("L73").Value = 6396
$ws.Range("M73").Value = -4570.2
$ws.Range("N73").Value = -8268
$ws.Range("H123").Value = 29999.666
$ws.Range("J123").Value = 29999.666
$ws.Range("L123").Value = 29999.666
$ws.Range("N123").Value = -34899.666

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 80670.46000000001
$ws.Range("I7").Value = 168801.17
$ws.Range("J7").Value = 5129.857
$ws.Range("K7").Value = 168801.17
$ws.Range("L7").Value = 5129.857
$ws.Range("M7").Value = -168689.17
$ws.Range("N7").Value = -5353.857
$ws.Range("H55").Value = 177.26666
$ws.Range("I55").Value = 116.75
$ws.Range("K55").Value = 116.75
$ws.Range("M55").Value = 56.25
$ws.Range("H126").Value = 80670.46000000001
$ws.Range("I126").Value = 168801.17
$ws.Range("J126").Value = 5129.857
$ws.Range("K126").Value = 506403.51
$ws.Range("L126").Value = 15389.571
$ws.Range("M126").Value = -503933.51
$ws.Range("N126").Value = -20329.571
$ws.Range("H132").Value = 4134.1
$ws.Range("J132").Value = 7799.6665
$ws.Range("L132").Value = 23398.9995
$ws.Range("N132").Value = -28458.9995

